# Apply the updated cryptocurrency price / volume data scraped by the
# GitHub Actions job to the worksheet, cell by cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new "Price" values look like plain numbers (e.g. "217.07")
# even though every cell in this sheet is stored as text in the source data
# (prices that do contain a thousands separator, like "26.861.58", are not
# ambiguous and do not need this treatment). Without help, Excel would silently
# reinterpret such a value as a floating point number when it is assigned,
# which would corrupt cells like "0.0629" or "217.07". Forcing the cell to a
# Text format before writing the value keeps it as the exact original string,
# and resetting the style back to Normal afterwards avoids leaving behind any
# formatting difference versus the rest of the sheet.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0629'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.26'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.530'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.85'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '215.07'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.41'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.49'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.37'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '147.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.73'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0508'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.37'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.02'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.810'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.24'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '61.99'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.08'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.66'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.408'
$ws.Range('D51').Style = 'Normal'

# The remaining updated cells (percentages and prices that already contain
# non-numeric characters) can be written directly.
$ws.Range('D2').Value = '26.861.58'
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').Value = '1.649.27'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('E4').Value = '  +0.68%  '
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('E7').Value = '  +0.53%  '
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('E9').Value = '  +0.28%  '
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('E11').Value = '  +0.35%  '
$ws.Range('D12').Value = '1.649.84'
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('E15').Value = '  -0.87%  '
$ws.Range('D16').Value = '26.840.63'
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('E17').Value = '  -1.02%  '
$ws.Range('E18').Value = '  -0.98%  '
$ws.Range('E19').Value = '  +0.58%  '
$ws.Range('E20').Value = '  +1.13%  '
$ws.Range('E21').Value = '  +11.23%  '
$ws.Range('E22').Value = '  -0.55%  '
$ws.Range('E23').Value = '  -1.37%  '
$ws.Range('E24').Value = '  +1.41%  '
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('E26').Value = '  -0.81%  '
$ws.Range('E27').Value = '  +0.48%  '
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('E31').Value = '  -0.40%  '
$ws.Range('E32').Value = '  -0.84%  '
$ws.Range('D33').Value = '1.297.08'
$ws.Range('E33').Value = '  +1.49%  '
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('E35').Value = '  +1.61%  '
$ws.Range('E36').Value = '  -1.98%  '
$ws.Range('E37').Value = '  +0.95%  '
$ws.Range('E38').Value = '  -0.58%  '
$ws.Range('E39').Value = '  +0.56%  '
$ws.Range('E40').Value = '  -0.64%  '
$ws.Range('E41').Value = '  -0.48%  '
$ws.Range('E42').Value = '  -2.14%  '
$ws.Range('D43').Value = '1.786.96'
$ws.Range('E43').Value = '  +0.26%  '
$ws.Range('E44').Value = '  +3.57%  '
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('E46').Value = '  +1.44%  '
$ws.Range('E47').Value = '  -0.93%  '
$ws.Range('E48').Value = '  +0.83%  '
$ws.Range('E49').Value = '  -1.98%  '
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('E51').Value = '  +0.39%  '
